$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.8
$ws.Range("G3").Value = 35
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 3.8
$ws.Range("G6").Value = 35
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 16.5
$ws.Range("G9").Value = 0.33

$ws.Range("G10").Select()
